$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.521.14'
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").Value = '  +0.48%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.727.49'
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").Value = '  +0.15%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9990'
$ws.Range("D4").NumberFormat = "General"
$ws.Range("E4").Value = '  +0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '245.29'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = '  +2.18%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9998'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = '  +0.04%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4809'
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = '  +1.67%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2669'
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = '  +1.13%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06220'
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = '  -0.36%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.725.57'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = '  +0.13%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07145'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = '  +0.86%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.67'
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = '  +1.74%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.6167'
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = '  +3.63%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.523'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = '  +2.56%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '77.13'
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = '  +0.89%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.9997'
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = '  +0.01%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.520.19'
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = '  +0.54%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.9995'
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = '  +0.03%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000006933'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = '  +1.68%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.67'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = '  +0.56%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.947.30'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = '  +0.50%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.528'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = '  -0.83%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.957'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = '  +1.95%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.288'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = '  -1.11%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '136.51'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = '  +1.27%  '

$ws.Range("E26").Value = '  +0.43%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.793'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = '  +1.17%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.403'
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = '  -0.57%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '106.79'
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = '  -1.77%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.980'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = '  -1.22%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08029'
$ws.Range("D31").NumberFormat = "General"

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.708'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = '  +0.06%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04565'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = '  +2.08%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.9990'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = '  +0.01%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.614'
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = '  +0.10%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6378'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = '  +2.43%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9912'
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = '  +1.16%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.9282'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = '  -0.04%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.092'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = '  +9.07%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.416'
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = '  -0.17%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '104.80'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = '  -9.94%  '

$ws.Range("E42").Value = '  +0.55%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.01504'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = '  +1.54%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.614'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = '  +4.41%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.3907'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = '  +1.91%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '6.905'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = '  +9.73%  '

$ws.Range("E47").Value = '  +1.45%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.05330'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = '  +0.71%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '30.91'
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = '  +0.92%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.845'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = '  +1.88%  '

$ws.Range("E51").Value = '  +3.87%  '
